$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Kool & The Gang" / "Steppin' Into Love" / "happy" row (row 7, cols A:C)
$ws.Range("A7:C7").ClearContents()

# Update the view selection to match the new state
$ws.Range("C11").Select()
